$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.015.18'
$ws.Range("E2").Value = '  +5.69%  '
$ws.Range("D3").Value = '2.418.89'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").Value = '''572.52'
$ws.Range("E5").Value = '  +2.44%  '
$ws.Range("D6").Value = '''145.96'
$ws.Range("E6").Value = '  +6.28%  '
$ws.Range("D7").Value = '''0.997'
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  +2.60%  '
$ws.Range("D9").Value = '2.451.18'
$ws.Range("E9").Value = '  +3.79%  '
$ws.Range("E10").Value = '  +6.02%  '
$ws.Range("D11").Value = '''0.160'
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("E12").Value = '  +3.11%  '
$ws.Range("E13").Value = '  +5.05%  '
$ws.Range("E14").Value = '  +7.36%  '
$ws.Range("D15").Value = '''0.0000177'
$ws.Range("E15").Value = '  +7.98%  '
$ws.Range("D16").Value = '2.861.08'
$ws.Range("E16").Value = '  +2.39%  '
$ws.Range("D17").Value = '62.946.62'
$ws.Range("E17").Value = '  +5.44%  '
$ws.Range("D18").Value = '2.458.23'
$ws.Range("E18").Value = '  +4.05%  '
$ws.Range("D19").Value = '''7.89'
$ws.Range("E19").Value = '  -1.60%  '
$ws.Range("D20").Value = '''10.98'
$ws.Range("E20").Value = '  +5.26%  '
$ws.Range("D21").Value = '''328.66'
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("E22").Value = '  +2.34%  '
$ws.Range("E23").Value = '  +12.88%  '
$ws.Range("D24").Value = '''0.998'
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("D25").Value = '''65.70'
$ws.Range("E25").Value = '  +2.59%  '
$ws.Range("D26").Value = '''640.08'
$ws.Range("E26").Value = '  +14.47%  '
$ws.Range("E27").Value = '  +10.31%  '
$ws.Range("D28").Value = '''8.50'
$ws.Range("E28").Value = '  +4.41%  '
$ws.Range("D29").Value = '0.0₃0987'
$ws.Range("E29").Value = '  +7.37%  '
$ws.Range("D30").Value = '2.532.99'
$ws.Range("D31").Value = '''8.20'
$ws.Range("E31").Value = '  +2.83%  '
$ws.Range("E32").Value = '  +9.04%  '
$ws.Range("E33").Value = '  +6.01%  '
$ws.Range("D34").Value = '''1.84'
$ws.Range("E34").Value = '  +3.90%  '
$ws.Range("E35").Value = '  +5.04%  '
$ws.Range("D36").Value = '''0.995'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("E37").Value = '  +5.07%  '
$ws.Range("E38").Value = '  +2.47%  '
$ws.Range("D39").Value = '''152.95'
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '''5.41'
$ws.Range("E40").Value = '  +8.92%  '
$ws.Range("D41").Value = '''18.69'
$ws.Range("E41").Value = '  +3.16%  '
$ws.Range("E42").Value = '  +14.19%  '
$ws.Range("D43").Value = '''1.77'
$ws.Range("E43").Value = '  +8.25%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '0.0₆0288'
$ws.Range("E45").Value = '  -3.29%  '
$ws.Range("D46").Value = '''144.85'
$ws.Range("E46").Value = '  +4.88%  '
$ws.Range("D47").Value = '''3.60'
$ws.Range("E47").Value = '  +2.51%  '
$ws.Range("D48").Value = '''20.41'
$ws.Range("E48").Value = '  +7.27%  '
$ws.Range("E49").Value = '  +3.37%  '
$ws.Range("D50").Value = '''0.0516'
$ws.Range("E50").Value = '  +3.52%  '
$ws.Range("D51").Value = '''12.52'
$ws.Range("E51").Value = '  +7.21%  '
